$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F5").Value = 64
$ws.Range("F6").Value = 815
$ws.Range("F7").Value = 391
$ws.Range("F8").Value = 4640
$ws.Range("F9").Value = 4640
$ws.Range("F11").Value = 117
$ws.Range("F12").Value = 151
$ws.Range("F15").Value = 106
$ws.Range("F16").Value = 7326
$ws.Range("F17").Value = 247
$ws.Range("F22").Value = 1332
$ws.Range("G22").Value = 63
$ws.Range("F25").Value = 1732
$ws.Range("F27").Value = 1988
$ws.Range("F28").Value = 6137
$ws.Range("F34").Value = 6345
$ws.Range("F40").Value = 12
$ws.Range("F41").Value = 2444
$ws.Range("F45").Value = 34
$ws.Range("F46").Value = 412
$ws.Range("C47").Value = "北京·第二届城市梦想动漫嘉年华（CDS）"
$ws.Range("F47").Value = 2122
$ws.Range("C49").Value = "北京·万游引力国潮动漫嘉年华s7"
$ws.Range("F49").Value = 1067

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F3").Value = 226
$ws.Range("F5").Value = 43
$ws.Range("F6").Value = 117

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F5").Value = 226
$ws.Range("F6").Value = 64
$ws.Range("F7").Value = 43
$ws.Range("F8").Value = 391
$ws.Range("F9").Value = 4640
$ws.Range("F10").Value = 4640
$ws.Range("F12").Value = 117
$ws.Range("F13").Value = 151
$ws.Range("F16").Value = 106
$ws.Range("F17").Value = 7326
$ws.Range("F18").Value = 247
$ws.Range("F21").Value = 1332
$ws.Range("G21").Value = 63
$ws.Range("F22").Value = 117
$ws.Range("F24").Value = 1732
$ws.Range("F26").Value = 1988
$ws.Range("F29").Value = 6137
$ws.Range("F36").Value = 6345
$ws.Range("F42").Value = 2444
$ws.Range("F45").Value = 34
$ws.Range("F46").Value = 412
$ws.Range("C48").Value = "北京·第二届城市梦想动漫嘉年华（CDS）"
$ws.Range("F48").Value = 2122
